$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting A:D to B:E
$ws.Columns.Item(1).Insert()

# New header cell
$ws.Range("A1").Value = "ID"

# Copy formatting (bold/border/centered) from the adjacent header cell
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the ID column with the row labels
$ws.Cells.Item(2, 1).Value = '35 B2Op'
$ws.Cells.Item(3, 1).Value = '36 B3Pop'
$ws.Cells.Item(4, 1).Value = '37 B4Pop'
$ws.Cells.Item(5, 1).Value = '38 B5C'
$ws.Cells.Item(6, 1).Value = '1 D1V'
$ws.Cells.Item(7, 1).Value = '3 D3V'
$ws.Cells.Item(8, 1).Value = '4 D4V'
$ws.Cells.Item(9, 1).Value = '5 D5V'
$ws.Cells.Item(10, 1).Value = '40 D2V'
$ws.Cells.Item(11, 1).Value = '42 D4V'
$ws.Cells.Item(12, 1).Value = '6 FB1C'
$ws.Cells.Item(13, 1).Value = '7 FB2C'
$ws.Cells.Item(14, 1).Value = '8 FB3C'
$ws.Cells.Item(15, 1).Value = '9 FB4C'
$ws.Cells.Item(16, 1).Value = '44 FB1C'
$ws.Cells.Item(17, 1).Value = '45 FB2C'
$ws.Cells.Item(18, 1).Value = '46 FB3Op'
$ws.Cells.Item(19, 1).Value = '47 FB4Pop'
$ws.Cells.Item(20, 1).Value = '48 FB5Pop'
$ws.Cells.Item(21, 1).Value = '10 H1C'
$ws.Cells.Item(22, 1).Value = '11 H2C'
$ws.Cells.Item(23, 1).Value = '12 H3C'
$ws.Cells.Item(24, 1).Value = '13 H4C'
$ws.Cells.Item(25, 1).Value = '14 H5C'
$ws.Cells.Item(26, 1).Value = '49 H1C'
$ws.Cells.Item(27, 1).Value = '51 H3C'
$ws.Cells.Item(28, 1).Value = '52 H4C'
$ws.Cells.Item(29, 1).Value = '53 H5De'
$ws.Cells.Item(30, 1).Value = '54 HH1De'
$ws.Cells.Item(31, 1).Value = '58 HH5De'
$ws.Cells.Item(32, 1).Value = '59 SF1C'
$ws.Cells.Item(33, 1).Value = '60 SF2C'
$ws.Cells.Item(34, 1).Value = '61 SF3C'
$ws.Cells.Item(35, 1).Value = '63 SF5C'
$ws.Cells.Item(36, 1).Value = '64 SLOp'
$ws.Cells.Item(37, 1).Value = '65 SOC'
$ws.Cells.Item(38, 1).Value = '66 ST1Rü'
$ws.Cells.Item(39, 1).Value = '67 ST2Rü'
$ws.Cells.Item(40, 1).Value = '68 ST3Rü'
$ws.Cells.Item(41, 1).Value = '69 ST4Rü'
$ws.Cells.Item(42, 1).Value = '67 ST5Rü'
$ws.Cells.Item(43, 1).Value = '15 Z1C'
$ws.Cells.Item(44, 1).Value = '16 Z2C'
$ws.Cells.Item(45, 1).Value = '17 Z3C'
$ws.Cells.Item(46, 1).Value = '18 Z4C'
$ws.Cells.Item(47, 1).Value = '19 Z5C'
